$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Extend the Task4 header text (F1 / shared string "Task₄")
$ws.Range("F1").Value = "Task₄ parcial 2 del 12 de octubre derivacion e integracion"

# 2. Fill in grades for row 5 (camila): D5 and E5 = 5
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 5

# 3. Fill in grade for row 9 (duvan salgado): E9 = 4.5
$ws.Range("E9").Value = 4.5

# 4. Add the underline-style marker cells below the table (rows 10 and 13, column E)
$ws.Range("E10").Font.Underline = 2
$ws.Range("E13").Font.Underline = 2

# 5. Widen column E to fit the new, longer header text
$ws.Columns("E").ColumnWidth = 17.75

# 6. Update the active selection to D6
[void]$ws.Range("D6").Select()
